$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 54.69462833333333
$ws.Cells.Item(2, 8).Value = 164.083885
$ws.Cells.Item(2, 9).Value = 0.2790924419198448
$ws.Cells.Item(2, 10).Value = 0.2790924419198448
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.1591403333333333
$ws.Cells.Item(2, 14).Value = 0.477421
$ws.Cells.Item(2, 15).Value = 0.01260326986877891
$ws.Cells.Item(2, 16).Value = 0.01260326986877891
$ws.Cells.Item(2, 17).Value = 8.704121384509444
$ws.Cells.Item(2, 18).Value = 78.337092460585
$ws.Cells.Item(2, 19).Value = 0.003517477363852307
$ws.Cells.Item(2, 20).Value = 0.003517477363852308

$ws.Cells.Item(3, 7).Value = 54.69462833333333
$ws.Cells.Item(3, 8).Value = 164.083885
$ws.Cells.Item(3, 9).Value = 0.2790924419198448
$ws.Cells.Item(3, 10).Value = 0.2790924419198448
$ws.Cells.Item(3, 15).Value = 0.3005157372251983
$ws.Cells.Item(3, 16).Value = 0.3005157372251983
$ws.Cells.Item(3, 17).Value = 207.5433980226989
$ws.Cells.Item(3, 18).Value = 1867.89058220429
$ws.Cells.Item(3, 19).Value = 0.08387167093752299
$ws.Cells.Item(3, 20).Value = 0.083871670937523

$ws.Cells.Item(4, 7).Value = 54.69462833333333
$ws.Cells.Item(4, 8).Value = 164.083885
$ws.Cells.Item(4, 9).Value = 0.2790924419198448
$ws.Cells.Item(4, 10).Value = 0.2790924419198448
$ws.Cells.Item(4, 13).Value = 8.673183333333334
$ws.Cells.Item(4, 14).Value = 26.01955
$ws.Cells.Item(4, 15).Value = 0.6868809929060228
$ws.Cells.Item(4, 16).Value = 0.6868809929060229
$ws.Cells.Item(4, 17).Value = 474.3765388835278
$ws.Cells.Item(4, 18).Value = 4269.388849951751
$ws.Cells.Item(4, 19).Value = 0.1917032936184695
$ws.Cells.Item(4, 20).Value = 0.1917032936184695

$ws.Cells.Item(5, 7).Value = 19.32115333333334
$ws.Cells.Item(5, 8).Value = 57.96346000000001
$ws.Cells.Item(5, 9).Value = 0.09859081282432611
$ws.Cells.Item(5, 10).Value = 0.09859081282432611
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.1591403333333333
$ws.Cells.Item(5, 14).Value = 0.477421
$ws.Cells.Item(5, 15).Value = 0.01260326986877891
$ws.Cells.Item(5, 16).Value = 0.01260326986877891
$ws.Cells.Item(5, 17).Value = 3.074774781851112
$ws.Cells.Item(5, 18).Value = 27.67297303666
$ws.Cells.Item(5, 19).Value = 0.00124256662060725
$ws.Cells.Item(5, 20).Value = 0.001242566620607251

$ws.Cells.Item(6, 7).Value = 19.32115333333334
$ws.Cells.Item(6, 8).Value = 57.96346000000001
$ws.Cells.Item(6, 9).Value = 0.09859081282432611
$ws.Cells.Item(6, 10).Value = 0.09859081282432611
$ws.Cells.Item(6, 15).Value = 0.3005157372251983
$ws.Cells.Item(6, 16).Value = 0.3005157372251983
$ws.Cells.Item(6, 17).Value = 73.31575218098224
$ws.Cells.Item(6, 18).Value = 659.8417696288401
$ws.Cells.Item(6, 19).Value = 0.02962809079953389
$ws.Cells.Item(6, 20).Value = 0.0296280907995339

$ws.Cells.Item(7, 7).Value = 19.32115333333334
$ws.Cells.Item(7, 8).Value = 57.96346000000001
$ws.Cells.Item(7, 9).Value = 0.09859081282432611
$ws.Cells.Item(7, 10).Value = 0.09859081282432611
$ws.Cells.Item(7, 13).Value = 8.673183333333334
$ws.Cells.Item(7, 14).Value = 26.01955
$ws.Cells.Item(7, 15).Value = 0.6868809929060228
$ws.Cells.Item(7, 16).Value = 0.6868809929060229
$ws.Cells.Item(7, 17).Value = 167.5759050714445
$ws.Cells.Item(7, 18).Value = 1508.183145643
$ws.Cells.Item(7, 19).Value = 0.06772015540418497
$ws.Cells.Item(7, 20).Value = 0.06772015540418497

$ws.Cells.Item(8, 7).Value = 11.023718
$ws.Cells.Item(8, 8).Value = 33.071154
$ws.Cells.Item(8, 9).Value = 0.05625116157486912
$ws.Cells.Item(8, 10).Value = 0.05625116157486911
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.1591403333333333
$ws.Cells.Item(8, 14).Value = 0.477421
$ws.Cells.Item(8, 15).Value = 0.01260326986877891
$ws.Cells.Item(8, 16).Value = 0.01260326986877891
$ws.Cells.Item(8, 17).Value = 1.754318157092667
$ws.Cells.Item(8, 18).Value = 15.788863413834
$ws.Cells.Item(8, 19).Value = 0.0007089485697603619
$ws.Cells.Item(8, 20).Value = 0.0007089485697603619

$ws.Cells.Item(9, 7).Value = 11.023718
$ws.Cells.Item(9, 8).Value = 33.071154
$ws.Cells.Item(9, 9).Value = 0.05625116157486912
$ws.Cells.Item(9, 10).Value = 0.05625116157486911
$ws.Cells.Item(9, 15).Value = 0.3005157372251983
$ws.Cells.Item(9, 16).Value = 0.3005157372251983
$ws.Cells.Item(9, 17).Value = 41.83043129245733
$ws.Cells.Item(9, 18).Value = 376.473881632116
$ws.Cells.Item(9, 19).Value = 0.01690435929044554
$ws.Cells.Item(9, 20).Value = 0.01690435929044554

$ws.Cells.Item(10, 7).Value = 11.023718
$ws.Cells.Item(10, 8).Value = 33.071154
$ws.Cells.Item(10, 9).Value = 0.05625116157486912
$ws.Cells.Item(10, 10).Value = 0.05625116157486911
$ws.Cells.Item(10, 13).Value = 8.673183333333334
$ws.Cells.Item(10, 14).Value = 26.01955
$ws.Cells.Item(10, 15).Value = 0.6868809929060228
$ws.Cells.Item(10, 16).Value = 0.6868809929060229
$ws.Cells.Item(10, 17).Value = 95.61072722896668
$ws.Cells.Item(10, 18).Value = 860.4965450607001
$ws.Cells.Item(10, 19).Value = 0.03863785371466322
$ws.Cells.Item(10, 20).Value = 0.03863785371466322

$ws.Cells.Item(11, 7).Value = 110.9336623333333
$ws.Cells.Item(11, 8).Value = 332.800987
$ws.Cells.Item(11, 9).Value = 0.5660655836809599
$ws.Cells.Item(11, 10).Value = 0.5660655836809599
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 0.6666666666666666
$ws.Cells.Item(11, 13).Value = 0.1591403333333333
$ws.Cells.Item(11, 14).Value = 0.477421
$ws.Cells.Item(11, 15).Value = 0.01260326986877891
$ws.Cells.Item(11, 16).Value = 0.01260326986877891
$ws.Cells.Item(11, 17).Value = 17.65402000161411
$ws.Cells.Item(11, 18).Value = 158.886180014527
$ws.Cells.Item(11, 19).Value = 0.007134277314558988
$ws.Cells.Item(11, 20).Value = 0.00713427731455899

$ws.Cells.Item(12, 7).Value = 110.9336623333333
$ws.Cells.Item(12, 8).Value = 332.800987
$ws.Cells.Item(12, 9).Value = 0.5660655836809599
$ws.Cells.Item(12, 10).Value = 0.5660655836809599
$ws.Cells.Item(12, 15).Value = 0.3005157372251983
$ws.Cells.Item(12, 16).Value = 0.3005157372251983
$ws.Cells.Item(12, 17).Value = 420.9471741072442
$ws.Cells.Item(12, 18).Value = 3788.524566965198
$ws.Cells.Item(12, 19).Value = 0.1701116161976958
$ws.Cells.Item(12, 20).Value = 0.1701116161976959

$ws.Cells.Item(13, 7).Value = 110.9336623333333
$ws.Cells.Item(13, 8).Value = 332.800987
$ws.Cells.Item(13, 9).Value = 0.5660655836809599
$ws.Cells.Item(13, 10).Value = 0.5660655836809599
$ws.Cells.Item(13, 13).Value = 8.673183333333334
$ws.Cells.Item(13, 14).Value = 26.01955
$ws.Cells.Item(13, 15).Value = 0.6868809929060228
$ws.Cells.Item(13, 16).Value = 0.6868809929060229
$ws.Cells.Item(13, 17).Value = 962.1479912550944
$ws.Cells.Item(13, 18).Value = 8659.33192129585
$ws.Cells.Item(13, 19).Value = 0.3888196901687051
$ws.Cells.Item(13, 20).Value = 0.3888196901687052

$wb.Save()
